$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the "level ideas" block (old rows 46-49) down by one row, to make
#    room for two new rows (44 and 45) describing the new subroutine items.
#    Work bottom-up so we never clobber a source cell before it's copied.
#    (Formats and values are copied separately: PasteSpecial xlPasteAll does
#    not reproduce the style reliably in this engine, so copy the format via
#    xlPasteFormats and then set the text explicitly.)
# ---------------------------------------------------------------------------

# old row49 -> new row50  ("everything you know is wrong")
$v = $ws.Cells.Item(49,3).Value2
$ws.Cells.Item(49,3).Copy()
$ws.Cells.Item(50,3).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(50,3).Value2 = $v

# old row48 -> new row49  ("bonus stages - short out levels ...")
$v = $ws.Cells.Item(48,3).Value2
$ws.Cells.Item(48,3).Copy()
$ws.Cells.Item(49,3).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(49,3).Value2 = $v

# old row47 -> new row48  ("ability to disable/enable end square")
$v = $ws.Cells.Item(47,3).Value2
$ws.Cells.Item(47,3).Copy()
$ws.Cells.Item(48,3).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(48,3).Value2 = $v

# old row46 -> new row47  ("level ideas" section header)
$v = $ws.Cells.Item(46,3).Value2
$ws.Cells.Item(46,3).Copy()
$ws.Cells.Item(47,3).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(47,3).Value2 = $v

# row 46 becomes an empty spacer row again (like rows 22/23/38 elsewhere)
$ws.Cells.Item(46,3).Clear()

# ---------------------------------------------------------------------------
# 2. Populate the two new rows (44 and 45) with the new wishlist items,
#    reusing formatting from similar existing rows so styles line up.
# ---------------------------------------------------------------------------

# Row 44: new bug/task note, formatted like row 43 (plain "item" row, col C only)
$ws.Cells.Item(43,3).Copy()
$ws.Cells.Item(44,3).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(44,3).Value2 = "subs calling subs may not be working right still, not 100% sure"

# Row 45: new task assigned to Tom, formatted like row 7 (A+C filled item row)
$ws.Cells.Item(7,1).Copy()
$ws.Cells.Item(45,1).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(45,1).Value2 = "Tom"

$ws.Cells.Item(7,3).Copy()
$ws.Cells.Item(45,3).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(45,3).Value2 = "game saves the level you're on and not the highest level you can select"

# ---------------------------------------------------------------------------
# 3. Update the view: scrolled position and active selection moved to A46.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A46").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
